# Daily BP terminal gate pricing (TGP) refresh.
# Shifts the two-day rolling window forward by one day: the newest date
# (effective date column A) moves from 18 Nov 2025 (45979) to 19 Nov 2025
# (45980), the previous 18 Nov values slide down into the 15 Nov rows (45976
# -> 45979), and fresh Diesel/ULP/PULP/e10 prices are written for the new day.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").Value2 = 45980
$ws.Range("D8").Value2 = 171.9
$ws.Range("E8").Value2 = 161.93
$ws.Range("F8").Value2 = 171.93
$ws.Range("G8").Value2 = 162.09

# Row 9
$ws.Range("A9").Value2 = 45980
$ws.Range("D9").Value2 = 171.9
$ws.Range("E9").Value2 = 161.93
$ws.Range("F9").Value2 = 171.93
$ws.Range("G9").Value2 = 162.09

# Row 10
$ws.Range("A10").Value2 = 45980
$ws.Range("D10").Value2 = 174
$ws.Range("E10").Value2 = 164.99
$ws.Range("F10").Value2 = 174.99
$ws.Range("G10").Value2 = 165.49

# Row 11
$ws.Range("A11").Value2 = 45979
$ws.Range("D11").Value2 = 171.98
$ws.Range("E11").Value2 = 161.81
$ws.Range("F11").Value2 = 171.81
$ws.Range("G11").Value2 = 161.97999999999999

# Row 12
$ws.Range("A12").Value2 = 45979
$ws.Range("D12").Value2 = 171.98
$ws.Range("E12").Value2 = 161.81
$ws.Range("F12").Value2 = 171.81
$ws.Range("G12").Value2 = 161.97999999999999

# Row 13
$ws.Range("A13").Value2 = 45979
$ws.Range("D13").Value2 = 174.07
$ws.Range("E13").Value2 = 164.86
$ws.Range("F13").Value2 = 174.86
$ws.Range("G13").Value2 = 165.36

# Row 17
$ws.Range("A17").Value2 = 45980
$ws.Range("D17").Value2 = 177.33
$ws.Range("E17").Value2 = 167.72
$ws.Range("F17").Value2 = 177.72

# Row 18
$ws.Range("A18").Value2 = 45979
$ws.Range("D18").Value2 = 177.4
$ws.Range("E18").Value2 = 167.62
$ws.Range("F18").Value2 = 177.62

# Row 22
$ws.Range("A22").Value2 = 45980
$ws.Range("D22").Value2 = 172.92
$ws.Range("E22").Value2 = 163.96
$ws.Range("F22").Value2 = 173.56
$ws.Range("G22").Value2 = 165.24

# Row 23
$ws.Range("A23").Value2 = 45980
$ws.Range("D23").Value2 = 178.78
$ws.Range("E23").Value2 = 168.49
$ws.Range("F23").Value2 = 178.49

# Row 24
$ws.Range("A24").Value2 = 45980
$ws.Range("D24").Value2 = 178.58
$ws.Range("E24").Value2 = 168.73
$ws.Range("F24").Value2 = 178.73

# Row 25
$ws.Range("A25").Value2 = 45980
$ws.Range("D25").Value2 = 179.41
$ws.Range("E25").Value2 = 168.14
$ws.Range("F25").Value2 = 178.14
$ws.Range("G25").Value2 = 168.18

# Row 26
$ws.Range("A26").Value2 = 45980
$ws.Range("D26").Value2 = 178.1
$ws.Range("E26").Value2 = 169.72
$ws.Range("F26").Value2 = 179.72

# Row 27
$ws.Range("A27").Value2 = 45979
$ws.Range("D27").Value2 = 172.99
$ws.Range("E27").Value2 = 163.83000000000001
$ws.Range("F27").Value2 = 173.43
$ws.Range("G27").Value2 = 165.12

# Row 28
$ws.Range("A28").Value2 = 45979
$ws.Range("D28").Value2 = 178.85
$ws.Range("E28").Value2 = 168.36
$ws.Range("F28").Value2 = 178.36

# Row 29
$ws.Range("A29").Value2 = 45979
$ws.Range("D29").Value2 = 178.65
$ws.Range("E29").Value2 = 168.6
$ws.Range("F29").Value2 = 178.6

# Row 30
$ws.Range("A30").Value2 = 45979
$ws.Range("D30").Value2 = 179.47
$ws.Range("E30").Value2 = 168.01
$ws.Range("F30").Value2 = 178.01
$ws.Range("G30").Value2 = 168.05

# Row 31
$ws.Range("A31").Value2 = 45979
$ws.Range("D31").Value2 = 178.17
$ws.Range("E31").Value2 = 169.58
$ws.Range("F31").Value2 = 179.58

# Row 35
$ws.Range("A35").Value2 = 45980
$ws.Range("D35").Value2 = 172.36
$ws.Range("E35").Value2 = 161.99
$ws.Range("F35").Value2 = 170.99

# Row 36
$ws.Range("A36").Value2 = 45979
$ws.Range("D36").Value2 = 172.54
$ws.Range("E36").Value2 = 161.86000000000001
$ws.Range("F36").Value2 = 170.86

# Row 40
$ws.Range("A40").Value2 = 45980
$ws.Range("D40").Value2 = 177.99
$ws.Range("E40").Value2 = 167.64
$ws.Range("F40").Value2 = 177.64

# Row 41
$ws.Range("A41").Value2 = 45980
$ws.Range("D41").Value2 = 177.69
$ws.Range("E41").Value2 = 168.06
$ws.Range("F41").Value2 = 178.06

# Row 42
$ws.Range("A42").Value2 = 45979
$ws.Range("D42").Value2 = 178.05
$ws.Range("E42").Value2 = 167.54
$ws.Range("F42").Value2 = 177.54

# Row 43
$ws.Range("A43").Value2 = 45979
$ws.Range("D43").Value2 = 177.75
$ws.Range("E43").Value2 = 167.95
$ws.Range("F43").Value2 = 177.95

# Row 47
$ws.Range("A47").Value2 = 45980
$ws.Range("D47").Value2 = 172.7
$ws.Range("E47").Value2 = 163.41999999999999
$ws.Range("F47").Value2 = 173.42

# Row 48
$ws.Range("A48").Value2 = 45980
$ws.Range("D48").Value2 = 172.66
$ws.Range("E48").Value2 = 163.58000000000001
$ws.Range("F48").Value2 = 173.58

# Row 49
$ws.Range("A49").Value2 = 45979
$ws.Range("D49").Value2 = 173.34
$ws.Range("E49").Value2 = 163.30000000000001
$ws.Range("F49").Value2 = 173.3

# Row 50
$ws.Range("A50").Value2 = 45979
$ws.Range("D50").Value2 = 173.31
$ws.Range("E50").Value2 = 163.46
$ws.Range("F50").Value2 = 173.46

# Row 54
$ws.Range("A54").Value2 = 45980
$ws.Range("D54").Value2 = 188.05
$ws.Range("E54").Value2 = 178.19
$ws.Range("F54").Value2 = 188.19

# Row 55
$ws.Range("A55").Value2 = 45980
$ws.Range("D55").Value2 = 175.75
$ws.Range("E55").Value2 = 175.14
$ws.Range("F55").Value2 = 185.14

# Row 56
$ws.Range("A56").Value2 = 45980
$ws.Range("D56").Value2 = 178.26

# Row 57
$ws.Range("A57").Value2 = 45980
$ws.Range("D57").Value2 = 177.76
$ws.Range("E57").Value2 = 169.41

# Row 58
$ws.Range("A58").Value2 = 45980
$ws.Range("D58").Value2 = 173.67
$ws.Range("E58").Value2 = 165.46
$ws.Range("F58").Value2 = 175.46

# Row 59
$ws.Range("A59").Value2 = 45980
$ws.Range("D59").Value2 = 180.29
$ws.Range("E59").Value2 = 176.19

# Row 60
$ws.Range("A60").Value2 = 45979
$ws.Range("D60").Value2 = 188.11
$ws.Range("E60").Value2 = 178.04
$ws.Range("F60").Value2 = 188.05

# Row 61
$ws.Range("A61").Value2 = 45979
$ws.Range("D61").Value2 = 175.82
$ws.Range("E61").Value2 = 175.02
$ws.Range("F61").Value2 = 185.02

# Row 62
$ws.Range("A62").Value2 = 45979
$ws.Range("D62").Value2 = 178.32

# Row 63
$ws.Range("A63").Value2 = 45979
$ws.Range("D63").Value2 = 177.83
$ws.Range("E63").Value2 = 169.29

# Row 64
$ws.Range("A64").Value2 = 45979
$ws.Range("D64").Value2 = 173.74
$ws.Range("E64").Value2 = 165.34
$ws.Range("F64").Value2 = 175.34

# Row 65
$ws.Range("A65").Value2 = 45979
$ws.Range("D65").Value2 = 180.37
$ws.Range("E65").Value2 = 176.05
